$d = $word.ActiveDocument
$table = $d.Tables(1)

# 1) Table cell: "positionUpdate" -> "position" (Nom de balise column, row 2 of table 1)
#    Scope the find to just this cell so it doesn't also hit the heading text below.
$cell1 = $table.Rows(2).Cells(1)
$cell1.Range.Find.Execute("positionUpdate", $true, $false, $false, $false, $false,
                           $true, 1, $false, "position", 2)

# 2) Heading: "Objet geoPos" -> "Objet positionUpdate"
$d.Content.Find.Execute("Objet geoPos", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Objet positionUpdate", 2)

# 3) Fill the previously-empty "Exemple" cell (last cell of row 2) with a single space.
$cell6 = $table.Rows(2).Cells(6)
$cell6.Range.Text = " "
